$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Revision_History: log a new revision row (TFS# 3877 - add jobcode WACQ13)
# -----------------------------------------------------------------------
$rh = $wb.Worksheets.Item("Revision_History")

# Copy the date formatting from the row above so the new date cell reuses
# the existing numFmt (m/d/yyyy) instead of creating a new style entry.
$rh.Cells.Item(2,2).Copy()
$rh.Cells.Item(5,2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rh.Cells.Item(5,1).Value = 3
$rh.Cells.Item(5,2).Value = 42633
$rh.Cells.Item(5,3).Value = "Susmitha Palacherla"
$rh.Cells.Item(5,4).Value = 3877
$rh.Cells.Item(5,5).Value = "Add jobcode WACQ13 in Role_access table(Mark Hackman and Scott Potter)"

# -----------------------------------------------------------------------
# AT_Role_Access: add the new WACQ13 ("Sr Specialist, Quality (CS)")
# jobcode as CoachingAdmin (role 101) and WarningAdmin (role 103), and
# update the trailing SQL-snippet comments in column K to match.
# -----------------------------------------------------------------------
$ra = $wb.Worksheets.Item("AT_Role_Access")

$ra.Cells.Item(9,1).Value = "WACQ13"
$ra.Cells.Item(9,2).Value = "Sr Specialist, Quality (CS)"
$ra.Cells.Item(9,3).Value = 101
$ra.Cells.Item(9,4).Value = "CoachingAdmin"
$ra.Cells.Item(9,5).Value = 0
$ra.Cells.Item(9,6).Value = 1

$ra.Cells.Item(10,1).Value = "WACQ13"
$ra.Cells.Item(10,2).Value = "Sr Specialist, Quality (CS)"
$ra.Cells.Item(10,3).Value = 103
$ra.Cells.Item(10,4).Value = "WarningAdmin"
$ra.Cells.Item(10,5).Value = 0
$ra.Cells.Item(10,6).Value = 1

$ra.Cells.Item(17,11).Value = "           ('WACQ13','Sr Specialist, Quality (CS)',103,'WarningAdmin',0,1)"
$ra.Cells.Item(15,11).Value = "           ('WISY14','Sr Analyst, Systems',103,'WarningAdmin',0,1),"
$ra.Cells.Item(16,11).Value = "      ('WACQ13','Sr Specialist, Quality (CS)',101,'CoachingAdmin',0,1),"

# Leave the same end-of-edit selection/active sheet state as the author:
# cursor parked a couple rows/columns past the last typed SQL comment on
# AT_Role_Access, but with Revision_History as the active (tabSelected) sheet.
[void]$ra.Activate()
[void]$ra.Range("M22").Select()

[void]$rh.Activate()
[void]$rh.Range("A2").Select()
